$wb = $excel.ActiveWorkbook

# The scenario is being edited on the "Edit Repayment Schedule" sheet: a new
# "waittopageload1" step (amount 2000) is inserted right after the
# "clickonEditRepaymentSchedule" row (old row 11, now row 12), pushing the
# subsequent approve/disburse rows down by one.
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row at position 11, shifting existing rows 11-17 down to 12-18.
$ws.Rows(11).Insert()

$ws.Range("A11").Value = "waittopageload1"
$ws.Range("B11").Value = 2000

# Match the look/format used by the other "amount" style cell (B3, 2000)
# rather than the formatting inherited from the row above during the insert.
$ws.Range("B3").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# This sheet becomes the active one, with the new row selected.
$ws.Activate()
$ws.Range("A11:B11").Select()
